$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab14")

$ws.Range("C97").Value = 4.5540752227489802
$ws.Range("D97").Value = 0.34843122615263999
$ws.Range("E97").Value = 0.68496586726261999
$ws.Range("F97").Value = 0.69022881640837996
$ws.Range("G97").Value = 0.04513061001825
$ws.Range("H97").Value = 0.70881008184872996

$ws.Range("C98").Value = 5.2005227009455401
$ws.Range("D98").Value = 0.34283769006530002
$ws.Range("E98").Value = 0.73093719780444999
$ws.Range("F98").Value = 0.81080003082751995
$ws.Range("G98").Value = 0.091898612910880001
$ws.Range("H98").Value = 0.7763326416413
